$d = $word.ActiveDocument

# Locate the paragraph that currently reads:
# "-nächster Schritt – schauen ob die includierungen RealTime sind"
$oldParaText = "-nächster Schritt – schauen ob die includierungen RealTime sind"
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -eq ($oldParaText + "`r")) {
        $target = $cand
    }
}

if ($target -eq $null) {
    throw "Could not find the paragraph containing the old 'nächster Schritt' text"
}

$full  = $target.Range
$start = $full.Start
$end   = $full.End - 1          # exclude the trailing paragraph mark

$oldLen = $end - $start

# Insert the three replacement runs (as a WordprocessingML fragment so they
# stay as distinct <w:r> elements instead of being merged into one run) right
# before the old text.
$insPoint = $d.Range($start, $start)
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:r><w:t>-</w:t></w:r>' + `
    '<w:r><w:t>auch mit Datenbankanbinddung getestet. Die kleinen Ansichten in der Hauptauswahl</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> werden so gemacht wie beim Testen.</w:t></w:r>' + `
    '</w:p>'
$insPoint.InsertXML($xml)

# Remove the original run text, which got pushed after our new runs.
$insertedTextLen = "-auch mit Datenbankanbinddung getestet. Die kleinen Ansichten in der Hauptauswahl werden so gemacht wie beim Testen.".Length
$oldStart = $start + $insertedTextLen
$oldEnd   = $oldStart + $oldLen
$oldRange = $d.Range($oldStart, $oldEnd)
$oldRange.Text = ""

# The _GoBack bookmark now sits right where the old text used to end (right
# before the paragraph mark). Split it into its own, new, empty paragraph
# and insert the new "-Teile..." paragraph in between.
$bm  = $d.Bookmarks.Item("_GoBack")
$pos = $bm.Start

$splitPoint = $d.Range($pos, $pos)
$splitPoint.InsertBefore("`r-Teile einer Seite dynamisch geladen mithilfe von JQuery`r")
